# Auto-applied update to cryptos list (price/volume refresh + row 37/38 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to store the value as TEXT (not auto-converted to a number),
    # matching the original inline-string cell content, then reset the style back
    # to Normal/default so no extraneous number-format/style is left on the cell.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '28.739.15'
$ws.Range("E2").Value = '  +2.44%  '

Set-TextValue "D3" '1.807.07'
$ws.Range("E3").Value = '  +0.28%  '

Set-TextValue "D4" '0.9999'
$ws.Range("E4").Value = '  -0.03%  '

Set-TextValue "D5" '316.17'
$ws.Range("E5").Value = '  -0.07%  '

Set-TextValue "D6" '1.000'
$ws.Range("E6").Value = '  -0.02%  '

Set-TextValue "D7" '0.5371'
$ws.Range("E7").Value = '  -3.26%  '

Set-TextValue "D8" '0.3799'
$ws.Range("E8").Value = '  -0.79%  '

Set-TextValue "D9" '0.07542'
$ws.Range("E9").Value = '  -0.99%  '

Set-TextValue "D10" '42.65'
$ws.Range("E10").Value = '  -1.15%  '

$ws.Range("E11").Value = '  -1.13%  '

Set-TextValue "D12" '0.9998'
$ws.Range("E12").Value = '  -0.03%  '

Set-TextValue "D13" '21.02'
$ws.Range("E13").Value = '  -0.84%  '

$ws.Range("E14").Value = '  -0.45%  '

Set-TextValue "D15" '7.392'
$ws.Range("E15").Value = '  +3.15%  '

Set-TextValue "D16" '1.802.16'
$ws.Range("E16").Value = '  +0.15%  '

Set-TextValue "D17" '90.60'
$ws.Range("E17").Value = '  -1.35%  '

Set-TextValue "D18" '0.00001067'
$ws.Range("E18").Value = '  -1.18%  '

Set-TextValue "D19" '0.06446'
$ws.Range("E19").Value = '  -0.77%  '

$ws.Range("E20").Value = '  +0.04%  '

Set-TextValue "D21" '17.28'
$ws.Range("E21").Value = '  +0.42%  '

Set-TextValue "D22" '5.926'
$ws.Range("E22").Value = '  -0.58%  '

Set-TextValue "D23" '28.716.15'
$ws.Range("E23").Value = '  +2.24%  '

Set-TextValue "D25" '2.103'
$ws.Range("E25").Value = '  +0.19%  '

Set-TextValue "D26" '160.75'
$ws.Range("E26").Value = '  +2.98%  '

Set-TextValue "D27" '20.52'
$ws.Range("E27").Value = '  -0.52%  '

Set-TextValue "D28" '2.373'
$ws.Range("E28").Value = '  +0.04%  '

Set-TextValue "D29" '2.007.90'
$ws.Range("E29").Value = '  -0.11%  '

Set-TextValue "D30" '123.27'
$ws.Range("E30").Value = '  +0.38%  '

Set-TextValue "D31" '1.111'
$ws.Range("E31").Value = '  -2.87%  '

Set-TextValue "D32" '0.1056'
$ws.Range("E32").Value = '  +1.24%  '

Set-TextValue "D33" '5.675'
$ws.Range("E33").Value = '  -0.84%  '

Set-TextValue "D34" '3.691'
$ws.Range("E34").Value = '  +1.87%  '

Set-TextValue "D35" '0.2273'
$ws.Range("E35").Value = '  +7.07%  '

Set-TextValue "D36" '0.06486'
$ws.Range("E36").Value = '  +7.37%  '

$ws.Range("E39").Value = '  +0.94%  '

$ws.Range("E40").Value = '  -1.94%  '

Set-TextValue "D41" '0.6273'
$ws.Range("E41").Value = '  -0.29%  '

Set-TextValue "D42" '1.208'
$ws.Range("E42").Value = '  +4.83%  '

Set-TextValue "D43" '0.9994'
$ws.Range("E43").Value = '  +0.00%  '

Set-TextValue "D44" '1.395'
$ws.Range("E44").Value = '  -0.47%  '

Set-TextValue "D45" '13.41'
$ws.Range("E45").Value = '  -0.19%  '

Set-TextValue "D46" '0.5896'
$ws.Range("E46").Value = '  -0.28%  '

Set-TextValue "D47" '3.688'
$ws.Range("E47").Value = '  +0.40%  '

$ws.Range("E48").Value = '  +3.66%  '

Set-TextValue "D49" '1.963'
$ws.Range("E49").Value = '  +1.85%  '

$ws.Range("E50").Value = '  +1.63%  '

Set-TextValue "D51" '0.06891'
$ws.Range("E51").Value = '  +1.42%  '

# Row 37 and 38 swapped rank position (coin identity + link + price + volume)
Set-TextValue "B37" 'VeChain'
Set-TextValue "C37" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D37" '0.02322'
$ws.Range("E37").Value = '  +0.76%  '

Set-TextValue "B38" 'FraxShare'
Set-TextValue "C38" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D38" '8.849'
$ws.Range("E38").Value = '  +2.49%  '
